{"js": "// Replace the division-problem text in each data cell of the table,\n// in document order, while preserving run/paragraph formatting.\n// table.getCell(row, col) gives us precise positional addressing so that\n// duplicate text (e.g. \"68\u00f79=\" appears twice in the source) is replaced\n// correctly at each position rather than ambiguously.\n\nconst replacements = [\n  // [rowIndex, colIndex, oldText, newText]\n  [0, 0, \"92\u00f72=\", \"14\u00f73=\"],\n  [0, 1, \"54\u00f75=\", \"85\u00f73=\"],\n  [0, 2, \"68\u00f79=\", \"42\u00f75=\"],\n  [0, 3, \"25\u00f79=\", \"88\u00f79=\"],\n  [0, 4, \"84\u00f73=\", \"41\u00f79=\"],\n\n  [4, 0, \"35\u00f74=\", \"82\u00f79=\"],\n  [4, 1, \"52\u00f79=\", \"60\u00f79=\"],\n  [4, 2, \"69\u00f77=\", \"95\u00f79=\"],\n  [4, 3, \"76\u00f76=\", \"57\u00f78=\"],\n  [4, 4, \"59\u00f72=\", \"95\u00f77=\"],\n\n  [8, 0, \"52\u00f75=\", \"21\u00f74=\"],\n  [8, 1, \"42\u00f78=\", \"97\u00f78=\"],\n  [8, 2, \"18\u00f75=\", \"19\u00f77=\"],\n  [8, 3, \"54\u00f77=\", \"26\u00f73=\"],\n  [8, 4, \"18\u00f77=\", \"16\u00f72=\"],\n\n  [12, 0, \"68\u00f76=\", \"36\u00f75=\"],\n  [12, 1, \"72\u00f76=\", \"29\u00f75=\"],\n  [12, 2, \"36\u00f73=\", \"59\u00f78=\"],\n  [12, 3, \"68\u00f79=\", \"46\u00f76=\"],\n  [12, 4, \"97\u00f75=\", \"76\u00f79=\"],\n\n  [16, 0, \"97\u00f76=\", \"99\u00f75=\"],\n  [16, 1, \"18\u00f76=\", \"85\u00f75=\"],\n  [16, 2, \"79\u00f73=\", \"36\u00f79=\"],\n  [16, 3, \"70\u00f74=\", \"41\u00f75=\"],\n  [16, 4, \"86\u00f79=\", \"96\u00f73=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst searchResults = [];\nfor (const [row, col, oldText] of replacements) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  searchResults.push(results);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, , , newText] = replacements[i];\n  const results = searchResults[i];\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the division-problem text in each data cell of the table,\n# in document order, while preserving run/paragraph formatting.\n#\n# Table.Cell(row, col) (1-indexed) gives precise positional addressing so\n# duplicate text (e.g. \"68\u00f79=\" appears twice in the source) is replaced\n# correctly at each position rather than ambiguously with a document-wide\n# Find/Replace.\n#\n# NOTE: we re-derive each cell's range via $d.Range($cell.Range.Start,\n# $cell.Range.End) (rather than using $cell.Range directly) so that\n# Find.Execute is correctly confined to that cell, and we pass\n# Replace:=wdReplaceOne (1) so only the single match in that scoped range\n# is substituted - this keeps duplicate source strings (like \"68\u00f79=\")\n# from all being replaced at once.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"92\u00f72=\"; New = \"14\u00f73=\" },\n    @{ Row = 1;  Col = 2; Old = \"54\u00f75=\"; New = \"85\u00f73=\" },\n    @{ Row = 1;  Col = 3; Old = \"68\u00f79=\"; New = \"42\u00f75=\" },\n    @{ Row = 1;  Col = 4; Old = \"25\u00f79=\"; New = \"88\u00f79=\" },\n    @{ Row = 1;  Col = 5; Old = \"84\u00f73=\"; New = \"41\u00f79=\" },\n\n    @{ Row = 5;  Col = 1; Old = \"35\u00f74=\"; New = \"82\u00f79=\" },\n    @{ Row = 5;  Col = 2; Old = \"52\u00f79=\"; New = \"60\u00f79=\" },\n    @{ Row = 5;  Col = 3; Old = \"69\u00f77=\"; New = \"95\u00f79=\" },\n    @{ Row = 5;  Col = 4; Old = \"76\u00f76=\"; New = \"57\u00f78=\" },\n    @{ Row = 5;  Col = 5; Old = \"59\u00f72=\"; New = \"95\u00f77=\" },\n\n    @{ Row = 9;  Col = 1; Old = \"52\u00f75=\"; New = \"21\u00f74=\" },\n    @{ Row = 9;  Col = 2; Old = \"42\u00f78=\"; New = \"97\u00f78=\" },\n    @{ Row = 9;  Col = 3; Old = \"18\u00f75=\"; New = \"19\u00f77=\" },\n    @{ Row = 9;  Col = 4; Old = \"54\u00f77=\"; New = \"26\u00f73=\" },\n    @{ Row = 9;  Col = 5; Old = \"18\u00f77=\"; New = \"16\u00f72=\" },\n\n    @{ Row = 13; Col = 1; Old = \"68\u00f76=\"; New = \"36\u00f75=\" },\n    @{ Row = 13; Col = 2; Old = \"72\u00f76=\"; New = \"29\u00f75=\" },\n    @{ Row = 13; Col = 3; Old = \"36\u00f73=\"; New = \"59\u00f78=\" },\n    @{ Row = 13; Col = 4; Old = \"68\u00f79=\"; New = \"46\u00f76=\" },\n    @{ Row = 13; Col = 5; Old = \"97\u00f75=\"; New = \"76\u00f79=\" },\n\n    @{ Row = 17; Col = 1; Old = \"97\u00f76=\"; New = \"99\u00f75=\" },\n    @{ Row = 17; Col = 2; Old = \"18\u00f76=\"; New = \"85\u00f75=\" },\n    @{ Row = 17; Col = 3; Old = \"79\u00f73=\"; New = \"36\u00f79=\" },\n    @{ Row = 17; Col = 4; Old = \"70\u00f74=\"; New = \"41\u00f75=\" },\n    @{ Row = 17; Col = 5; Old = \"86\u00f79=\"; New = \"96\u00f73=\" }\n)\n\nforeach ($item in $replacements) {\n    $cell = $tbl.Cell($item.Row, $item.Col)\n    $rng = $d.Range($cell.Range.Start, $cell.Range.End)\n    $rng.Find.Execute($item.Old, $false, $false, $false, $false, $false, $true, 1, $false, $item.New, 1)\n}\n"}
